$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit reflects a re-run of the data-imputation pipeline:
#  - two sample rows ("RM 232" and "SC 92") were dropped from the dataset,
#    shifting every following sample row up and leaving only 32 data rows
#    (header + 32 = A1:F33, down from A1:F35)
#  - a handful of individual cells switched between "imputed" (blank) and
#    "known" (numeric) values

# --- Rows 2-25 keep the same sample in the same row, but a few cells in
#     column D (and one in column E) toggle between blank and a value ---
$ws.Range("E2").ClearContents()
$ws.Range("D6").Value = -14.2
$ws.Range("D8").ClearContents()
$ws.Range("D18").Value = -15.2
$ws.Range("D20").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("D25").ClearContents()

# --- Rows 26-33: after removing "RM 232" (old row 26) and "SC 92" (old row
#     28), the remaining SC samples move up two rows; write their values in
#     directly (some also switch between imputed/known along the way) ---
$ws.Range("A26").Value = "SC 5"
$ws.Range("B26").Value = -20.2
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38

$ws.Range("A27").Value = "SC 101"
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17

$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").ClearContents()
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = 17.44

$ws.Range("A29").Value = "SC 119"
$ws.Range("B29").Value = -19.5
$ws.Range("C29").ClearContents()
$ws.Range("D29").Value = -13
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06

$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89

$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = -8.1
$ws.Range("F31").Value = 17.18

$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").ClearContents()
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = -6.4
$ws.Range("F32").Value = 17.39

$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53

# --- The table now ends at row 33; clear out what used to be rows 34-35 so
#     the sheet's used range / dimension shrinks from A1:F35 to A1:F33 ---
$ws.Range("A34:F35").ClearContents()
